# Apply updated optimality/feasibility cut values to master_solution_decisions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 10.48436093167934
$ws.Range("G2").Value = 117.6

$ws.Range("E3").Value = 1.929375
$ws.Range("G3").Value = 29.4

$ws.Range("F4").Value = 509.65772

$ws.Range("E14").Value = 19.37741327286923
$ws.Range("G14").Value = 1711.080000000003

$ws.Range("G15").Value = 52.9200000000003

$ws.Range("E16").Value = 19.79669867147875
$ws.Range("G16").Value = 1398.74

$ws.Range("E17").Value = 1.514100000000009
$ws.Range("G17").Value = 43.26000000000022
